$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "68.358.36"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "3.597.41"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'622.02"
$ws.Range("E5").Value = "  -7.32%  "
$ws.Range("D6").Value = "'156.10"
$ws.Range("E6").Value = "  -2.64%  "
$ws.Range("D7").Value = "3.593.66"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("D12").Value = "'0.435"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").Value = "4.210.18"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "'32.21"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "3.597.61"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "68.343.22"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'6.45"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'15.68"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").Value = "'459.70"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").Value = "'9.83"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'0.644"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "'78.13"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "3.742.59"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D27").Value = "'10.76"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  -8.35%  "
$ws.Range("E29").Value = "  -7.35%  "
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "'26.22"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("D35").Value = "3.598.92"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").Value = "'6.21"
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("D38").Value = "'8.15"
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'177.71"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E42").Value = "  -7.67%  "
$ws.Range("D43").Value = "'0.0884"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("E44").Value = "  -5.19%  "
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "'46.10"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").Value = "'28.36"
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("E50").Value = "  -6.44%  "
$ws.Range("E51").Value = "  -5.67%  "
